$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new time-log entry in row 20 (copy formats from the row above first
# so the new cells reuse the existing date/time/minutes/hours number styles)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A20").Value = 41884
$ws.Range("B20").Value = 0.42083333333333334
$ws.Range("D20").Value = 5
$ws.Range("C20").Value = 0.45833333333333331
$ws.Range("F20").Value = "Coding"

# Update selection to reflect where the user left off (A21)
$ws.Range("A21").Select()

$wb.Save()
